$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# New resazurin survival observations recorded 8/13 (immune treatment,
# bag 49) appended below the existing "data" sheet log as of 8/14.
$rows = @(
    @(20240813,49,"immune","18C",1,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",2,0,0,0,0,0,1),
    @(20240813,49,"immune","18C",3,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",4,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",5,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",6,0,0,0,0,0,1),
    @(20240813,49,"immune","18C",7,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",8,0,0,0,0,1,1),
    @(20240813,49,"immune","18C",9,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",10,0,0,0,0,0,1),
    @(20240813,49,"immune","18C",11,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",12,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",13,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",14,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",15,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",16,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",17,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",18,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",19,0,0,0,0,0,0),
    @(20240813,49,"immune","18C",20,0,0,0,0,0,1),
    @(20240813,49,"immune","42C",21,0,0,0,0,0,1),
    @(20240813,49,"immune","42C",22,0,0,1,1,1,1),
    @(20240813,49,"immune","42C",23,0,0,0,0,0,0),
    @(20240813,49,"immune","42C",24,0,0,0,0,0,0),
    @(20240813,49,"immune","42C",25,0,0,1,1,1,1),
    @(20240813,49,"immune","42C",26,0,0,0,0,0,0),
    @(20240813,49,"immune","42C",27,0,0,1,1,1,1),
    @(20240813,49,"immune","42C",28,0,0,1,1,1,1),
    @(20240813,49,"immune","42C",29,0,0,1,1,1,1),
    @(20240813,49,"immune","42C",30,0,0,1,1,1,1),
    @(20240813,49,"immune","42C",31,0,1,1,1,1,1),
    @(20240813,49,"immune","42C",32,0,0,1,1,1,1),
    @(20240813,49,"immune","42C",33,0,0,0,0,0,0),
    @(20240813,49,"immune","42C",34,0,1,1,1,1,1),
    @(20240813,49,"immune","42C",35,0,0,0,0,0,0),
    @(20240813,49,"immune","42C",36,0,1,1,1,1,1),
    @(20240813,49,"immune","42C",37,0,0,1,1,1,1),
    @(20240813,49,"immune","42C",38,0,0,0,0,0,0),
    @(20240813,49,"immune","42C",39,0,0,1,1,1,1),
    @(20240813,49,"immune","42C",40,0,0,1,1,1,1)
)

$startRow = 362
$numRows = $rows.Count
$numCols = 11

$arr = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}

$endRow = $startRow + $numRows - 1
$targetRange = $ws.Range($ws.Cells.Item($startRow,1), $ws.Cells.Item($endRow,11))
$targetRange.Value2 = $arr

# Move the sheet's active selection / scroll position down to the new
# last row of data, matching where the author left off entering values.
$ws.Range("K" + $endRow).Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = $startRow + 5
    $win.Panes.Item(2).ScrollRow = $startRow + 5
} catch {}

# Nudge the saved window position leftward, matching the author's
# workbook-level view metadata update.
try {
    $excel.ActiveWindow.Left = -620
} catch {}
